$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price data rows (2-10) got rotated: the values in columns
# D (Fecha), H (Variedad), I (Calidad), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg)
# were reassigned among the rows (rows 4 and 9 unchanged).

$ws.Range("D2").Value = 44414
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 6000
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = 6500
$ws.Range("P2").Value = 6500

$ws.Range("D3").Value = 44253
$ws.Range("H3").Value = "Americana (o)"
$ws.Range("I3").Value = "Segunda"
$ws.Range("K3").Value = 4000
$ws.Range("L3").Value = 4500
$ws.Range("M3").Value = 4250
$ws.Range("P3").Value = 4250

$ws.Range("D5").Value = 44309
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 8500
$ws.Range("P5").Value = 8500

$ws.Range("D6").Value = 44281
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 5500
$ws.Range("P6").Value = 5500

$ws.Range("D7").Value = 44371
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7375
$ws.Range("P7").Value = 7375

$ws.Range("D8").Value = 44263
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = 7500
$ws.Range("P8").Value = 7500

$ws.Range("D10").Value = 44497
$ws.Range("J10").Value = 160
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 6000
$ws.Range("M10").Value = 5500
$ws.Range("P10").Value = 5500
